# Zoo management workbook - "First version of reading the data from the file"
#
# Helper: some string values ("True" / "False") would be auto-coerced to
# native Excel booleans if assigned directly via .Value (same as typing them
# into a cell). The target file stores them as literal text (shared-string
# cells), so we route those specific assignments through a
# Formula("=""..." ") -> Copy -> PasteSpecial(xlPasteValues) round-trip,
# which keeps the text as a true string type.
function Set-TextValue {
    param($range, [string]$text)
    $escaped = $text.Replace("""", """""")
    $range.Formula = "=""" + $escaped + """"
    $range.Copy() | Out-Null
    $range.PasteSpecial(-4163) | Out-Null
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "users"
# ---------------------------------------------------------------------
$users = $wb.Worksheets.Item("users")

$users.Range("D2").ClearContents() | Out-Null

$users.Range("D3").Value = "1;4"
$users.Range("D4").Value = "3;"
$users.Range("D5").Value = "2;"

# ---------------------------------------------------------------------
# Sheet "sections" - values unchanged, nothing to do.
# ---------------------------------------------------------------------

# ---------------------------------------------------------------------
# Sheet "cages"
# ---------------------------------------------------------------------
$cages = $wb.Worksheets.Item("cages")

$cages.Range("B2").Value = "1;2;3"

$cages.Range("A3").Value = 2
$cages.Range("B3").Value = "4;6;"

$cages.Range("A4").Value = 3
$cages.Range("B4").Value = "5;"

$cages.Range("A5").Value = 4
$cages.Range("B5").Value = "7;8;"

$cages.Range("B6").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet "animals"
# ---------------------------------------------------------------------
$animals = $wb.Worksheets.Item("animals")

# Row 2 (id 1): Leo, mammal/Lion
$animals.Range("B2").Value = "Leo"
$animals.Range("C2").Value = 5
$animals.Range("D2").Value = "mammal"
$animals.Range("E2").Value = "Lion"
$animals.Range("F2").Value = "African Lion"
$animals.Range("G2").Value = "Golden"
$animals.Range("J2").ClearContents() | Out-Null

# Row 3 (id 2): Molly, mammal/Dog
$animals.Range("B3").Value = "Molly"
$animals.Range("C3").Value = 3
$animals.Range("D3").Value = "mammal"
$animals.Range("E3").Value = "Dog"
$animals.Range("F3").Value = "Labrador"
$animals.Range("G3").Value = "Black"
$animals.Range("J3").ClearContents() | Out-Null

# Row 4 (id 3): Snowy, mammal/Cat
$animals.Range("B4").Value = "Snowy"
$animals.Range("C4").Value = 2
$animals.Range("D4").Value = "mammal"
$animals.Range("E4").Value = "Cat"
$animals.Range("F4").Value = "British Shorthair"
$animals.Range("G4").Value = "White"
$animals.Range("J4").ClearContents() | Out-Null

# Row 5 (id 4): Sky, bird
$animals.Range("A5").Value = 4
$animals.Range("B5").Value = "Sky"
$animals.Range("C5").Value = 4
$animals.Range("D5").Value = "bird"
$animals.Range("H5").Value = 120
Set-TextValue $animals.Range("I5") "True"

# Row 6 (id 5): Polly, bird
$animals.Range("A6").Value = 5
$animals.Range("B6").Value = "Polly"
$animals.Range("C6").Value = 6
$animals.Range("D6").Value = "bird"
$animals.Range("H6").Value = 40
Set-TextValue $animals.Range("I6") "True"

# Row 7 (id 6): Oscar, bird
$animals.Range("A7").Value = 6
$animals.Range("B7").Value = "Oscar"
$animals.Range("C7").Value = 8
$animals.Range("D7").Value = "bird"
$animals.Range("H7").Value = 90
Set-TextValue $animals.Range("I7") "False"

# Row 8 (id 7): Rex, reptile
$animals.Range("A8").Value = 7
$animals.Range("B8").Value = "Rex"
$animals.Range("C8").Value = 7
$animals.Range("D8").Value = "reptile"
Set-TextValue $animals.Range("J8") "False"

# Row 9 (id 8): Slither, reptile
$animals.Range("A9").Value = 8
$animals.Range("B9").Value = "Slither"
$animals.Range("C9").Value = 4
$animals.Range("D9").Value = "reptile"
Set-TextValue $animals.Range("J9") "True"

# Row 10 (id 9): Spike, reptile
$animals.Range("A10").Value = 9
$animals.Range("B10").Value = "Spike"
$animals.Range("C10").Value = 10
$animals.Range("D10").Value = "reptile"
Set-TextValue $animals.Range("J10") "True"

$excel.CutCopyMode = 0

$animals.Range("I27").Select() | Out-Null

# Restore "users" as the active sheet (it was the tab selected originally)
# with the caret parked on its final selection cell. Must run last: each
# sheet remembers its own last selection independently, but only the sheet
# selected *last* keeps tabSelected="1" at save time.
$users.Range("H7").Select() | Out-Null
